$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.593.71"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "3.427.84"
$ws.Range("E3").Value = "  +0.98%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.35"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.96"
$ws.Range("E6").Value = "  -2.15%  "

$ws.Range("D7").Value = "3.421.64"
$ws.Range("E7").Value = "  +0.98%  "

$ws.Range("E9").Value = "  -0.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.201"
$ws.Range("E10").Value = "  +1.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.584"
$ws.Range("E11").Value = "  -1.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.78"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000282"
$ws.Range("E13").Value = "  -0.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "692.21"
$ws.Range("E14").Value = "  +1.19%  "

$ws.Range("D15").Value = "3.977.15"
$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("D17").Value = "69.643.98"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.427.72"
$ws.Range("E18").Value = "  +0.85%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.121"
$ws.Range("E19").Value = "  +1.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.68"
$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.36"
$ws.Range("E21").Value = "  -0.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.896"
$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.44"
$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.97"
$ws.Range("E24").Value = "  -1.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.27"
$ws.Range("E25").Value = "  -2.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.91"
$ws.Range("E26").Value = "  -0.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.66"
$ws.Range("E27").Value = "  -2.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.59"
$ws.Range("E28").Value = "  -0.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.53"
$ws.Range("E29").Value = "  -2.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.78"
$ws.Range("E30").Value = "  +0.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.09"
$ws.Range("E31").Value = "  +1.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "572.36"
$ws.Range("E32").Value = "  +3.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.70"
$ws.Range("E33").Value = "  +0.33%  "

$ws.Range("E34").Value = "  -1.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.23"
$ws.Range("E35").Value = "  +0.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.104"
$ws.Range("E36").Value = "  -2.76%  "

$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").Value = "3.574.34"
$ws.Range("E38").Value = "  -3.56%  "

$ws.Range("E39").Value = "  -0.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.95"
$ws.Range("E40").Value = "  -0.32%  "

$ws.Range("E41").Value = "  +3.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.26"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("E43").Value = "  -0.52%  "

$ws.Range("E45").Value = "  -2.18%  "

$ws.Range("E46").Value = "  -1.52%  "

$ws.Range("E47").Value = "  +4.66%  "

$ws.Range("E48").Value = "  -0.23%  "

$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.65"
$ws.Range("E51").Value = "  +0.44%  "
